# Auto-sync update 2026-01-12 14:05:07
# Insert a new user-code row ("RNL06004" / "WaLc-sBJY") so the list keeps
# its alphabetical ordering. The row belongs right before the existing
# "RSN71799" entry (currently row 143), so shift rows 143:184 down by one
# and populate the freshly opened row 143 with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(143).Insert()

$ws.Range("A143").Value = "RNL06004"
$ws.Range("B143").Value = "WaLc-sBJY"
